$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 41: continuation comment ("...and same will be done with owner click")
$ws.Range("C41").Value = "and same will be done with "

# New row 42: "owner click" comment continuation
$ws.Range("C42").Value = "owner click"

# Row 40 gains Owner/State-style columns + a comment describing the owner-click work
$ws.Range("D40").Value = "wip"
$ws.Range("E40").Value = "shilpa"
$ws.Range("G40").Value = "Owner Details,form Location using tab,Grid showof Form Details"

# New row 41 also gets a second comment column entry
$ws.Range("G41").Value = "Register User,Login Ticket,Password encoding,busy symbol"

# Move selection / view down to the newly added rows, like the author scrolled to row ~41
$ws.Activate() | Out-Null
$ws.Range("G41").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
